$d = $word.ActiveDocument

# 1. Date: "01" -> "25" (date line "01 July, 2019")
$d.Content.Find.Execute("01 July, 2019", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25 July, 2019", 2)

# 2. "surveys are regularly conducted to monitor" -> "surveys are conducted to monitor"
$d.Content.Find.Execute("surveys are regularly conducted to monitor", $true, $false, $false, $false, $false,
                         $true, 1, $false, "surveys are conducted to monitor", 2)

# 3. "new resistance genes need to be deployed." -> "new resistance genes are needed."
$d.Content.Find.Execute("new resistance genes need to be deployed.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "new resistance genes are needed.", 2)

# 4. "was used to obtain pathotype data." -> "was used for data analysis."
$d.Content.Find.Execute("was used to obtain pathotype data.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "was used for data analysis.", 2)

# 5. Remove "a single state's or multiple states'" before "P. sojae population(s)"
$d.Content.Find.Execute("in relation to a single state" + [char]8217 + "s or multiple states" + [char]8217, $true, $false, $false, $false, $false,
                         $true, 1, $false, "in relation to", 2)

# 6. "was produced to support" -> "was developed to support"
$d.Content.Find.Execute("was produced to support", $true, $false, $false, $false, $false,
                         $true, 1, $false, "was developed to support", 2)

# 7. "Michigan Soy," -> "Michigan Soybean Promotion Committee,"
$d.Content.Find.Execute("Michigan Soy,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Michigan Soybean Promotion Committee,", 2)
